# Applies the cell-value updates shown in the diff (Price/Volume(1h)/Hora
# columns, plus a few Coin/Link swaps caused by re-ranking of rows 8-17).
# Columns D, E and G hold numeric-looking data that is stored as *text* in
# the workbook (t="inlineStr"), so values are written with a leading
# apostrophe to force Excel to keep them as text instead of auto-converting
# them to numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'311.49"
$ws.Range("E2").Value = "'2.22%"
$ws.Range("G2").Value = "'17"
# Row 3
$ws.Range("D3").Value = "'35.36"
$ws.Range("E3").Value = "'-2.61%"
$ws.Range("G3").Value = "'17"
# Row 4
$ws.Range("D4").Value = "'5.084"
$ws.Range("E4").Value = "'1.18%"
$ws.Range("G4").Value = "'17"
# Row 5
$ws.Range("D5").Value = "'0.08152"
$ws.Range("E5").Value = "'3.56%"
$ws.Range("G5").Value = "'17"
# Row 6
$ws.Range("D6").Value = "'2.087"
$ws.Range("E6").Value = "'-2.32%"
$ws.Range("G6").Value = "'17"
# Row 7
$ws.Range("D7").Value = "'7.953"
$ws.Range("E7").Value = "'0.06%"
$ws.Range("G7").Value = "'17"
# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9293"
$ws.Range("E8").Value = "'0.90%"
$ws.Range("G8").Value = "'17"
# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1036"
$ws.Range("E9").Value = "'7.55%"
$ws.Range("G9").Value = "'17"
# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1909"
$ws.Range("E10").Value = "'2.91%"
$ws.Range("G10").Value = "'17"
# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09182"
$ws.Range("E11").Value = "'6.67%"
$ws.Range("G11").Value = "'17"
# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03654"
$ws.Range("E12").Value = "'2.06%"
$ws.Range("G12").Value = "'17"
# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09896"
$ws.Range("E13").Value = "'-0.31%"
$ws.Range("G13").Value = "'17"
# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001434"
$ws.Range("E14").Value = "'-0.22%"
$ws.Range("G14").Value = "'17"
# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005817"
$ws.Range("E15").Value = "'2.16%"
$ws.Range("G15").Value = "'17"
# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.467"
$ws.Range("E16").Value = "'-0.11%"
$ws.Range("G16").Value = "'17"
# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.135"
$ws.Range("E17").Value = "'-0.16%"
$ws.Range("G17").Value = "'17"
# Row 18
$ws.Range("D18").Value = "'2.975"
$ws.Range("E18").Value = "'8.10%"
$ws.Range("G18").Value = "'17"
# Row 19
$ws.Range("D19").Value = "'0.3411"
$ws.Range("E19").Value = "'1.07%"
$ws.Range("G19").Value = "'17"
# Row 20
$ws.Range("D20").Value = "'0.1300"
$ws.Range("E20").Value = "'-3.39%"
$ws.Range("G20").Value = "'17"
# Row 21
$ws.Range("D21").Value = "'5.118"
$ws.Range("E21").Value = "'-0.62%"
$ws.Range("G21").Value = "'17"
# Row 22
$ws.Range("E22").Value = "'-1.57%"
$ws.Range("G22").Value = "'17"
# Row 23
$ws.Range("D23").Value = "'0.04547"
$ws.Range("E23").Value = "'-0.49%"
$ws.Range("G23").Value = "'17"
# Row 24
$ws.Range("D24").Value = "'0.001230"
$ws.Range("E24").Value = "'-0.21%"
$ws.Range("G24").Value = "'17"
# Row 25
$ws.Range("D25").Value = "'0.004772"
$ws.Range("E25").Value = "'-0.52%"
$ws.Range("G25").Value = "'17"
# Row 26
$ws.Range("D26").Value = "'0.0001254"
$ws.Range("E26").Value = "'-3.62%"
$ws.Range("G26").Value = "'17"
# Row 27
$ws.Range("D27").Value = "'0.0004461"
$ws.Range("E27").Value = "'-6.13%"
$ws.Range("G27").Value = "'17"
# Row 28
$ws.Range("G28").Value = "'17"
# Row 29
$ws.Range("G29").Value = "'17"
# Row 30
$ws.Range("G30").Value = "'17"
# Row 31
$ws.Range("G31").Value = "'17"
# Row 32
$ws.Range("G32").Value = "'17"
# Row 33
$ws.Range("G33").Value = "'17"
# Row 34
$ws.Range("G34").Value = "'17"
# Row 35
$ws.Range("G35").Value = "'17"
# Row 36
$ws.Range("G36").Value = "'17"
# Row 37
$ws.Range("G37").Value = "'17"
# Row 38
$ws.Range("G38").Value = "'17"
# Row 39
$ws.Range("D39").Value = "'0.01974"
$ws.Range("E39").Value = "'6.86%"
$ws.Range("G39").Value = "'17"
# Row 40
$ws.Range("D40").Value = "'0.04901"
$ws.Range("E40").Value = "'3.94%"
$ws.Range("G40").Value = "'17"
# Row 41
$ws.Range("D41").Value = "'0.007536"
$ws.Range("E41").Value = "'-3.16%"
$ws.Range("G41").Value = "'17"
# Row 42
$ws.Range("D42").Value = "'0.1380"
$ws.Range("E42").Value = "'-0.24%"
$ws.Range("G42").Value = "'17"
# Row 43
$ws.Range("D43").Value = "'0.008181"
$ws.Range("E43").Value = "'5.58%"
$ws.Range("G43").Value = "'17"
# Row 44
$ws.Range("D44").Value = "'0.002170"
$ws.Range("E44").Value = "'-2.33%"
$ws.Range("G44").Value = "'17"
# Row 45
$ws.Range("D45").Value = "'0.01150"
$ws.Range("E45").Value = "'0.72%"
$ws.Range("G45").Value = "'17"
# Row 46
$ws.Range("D46").Value = "'0.00006591"
$ws.Range("E46").Value = "'3.31%"
$ws.Range("G46").Value = "'17"
# Row 47
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.23%"
$ws.Range("G47").Value = "'17"
# Row 48
$ws.Range("D48").Value = "'185.15"
$ws.Range("E48").Value = "'254.74%"
$ws.Range("G48").Value = "'17"
# Row 49
$ws.Range("D49").Value = "'0.001505"
$ws.Range("E49").Value = "'-20.88%"
$ws.Range("G49").Value = "'17"
# Row 50
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.23%"
$ws.Range("G50").Value = "'17"
# Row 51
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.23%"
$ws.Range("G51").Value = "'17"
